$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update freight (kg) column G values for the rows listed in the diff
$ws.Range("G7").Value = 250
$ws.Range("G8").Value = 200
$ws.Range("G10").Value = 250
$ws.Range("G11").Value = 200
$ws.Range("G12").Value = 200
$ws.Range("G14").Value = 130
$ws.Range("G15").Value = 130
$ws.Range("G18").Value = 250
$ws.Range("G19").Value = 100
$ws.Range("G20").Value = 100
$ws.Range("G22").Value = 250
$ws.Range("G24").Value = 250
$ws.Range("G25").Value = 100
$ws.Range("G27").Value = 100

# Remove the trailing rows describing the JAV <-> Swiss Camp (SWC) traverse
# and its "total fly time" summary row (rows 29, 30, 31)
$ws.Rows("29:31").Delete()
